$d = $word.ActiveDocument
$xml = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:pPr><w:rPr><w:rFonts w:ascii="Arial" w:hAnsi="Arial" w:cs="Arial"/><w:b/><w:color w:val="0D0D0D" w:themeColor="text1" w:themeTint="F2"/><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:ascii="Arial" w:hAnsi="Arial" w:cs="Arial"/><w:b/><w:color w:val="0D0D0D" w:themeColor="text1" w:themeTint="F2"/><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr><w:t>Objetivo:</w:t></w:r></w:p><w:p><w:pPr><w:rPr><w:rFonts w:ascii="Arial" w:hAnsi="Arial" w:cs="Arial"/><w:color w:val="0D0D0D" w:themeColor="text1" w:themeTint="F2"/><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr></w:pPr></w:p><w:p><w:pPr><w:rPr><w:rFonts w:ascii="Arial" w:hAnsi="Arial" w:cs="Arial"/><w:color w:val="0D0D0D" w:themeColor="text1" w:themeTint="F2"/><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:ascii="Arial" w:hAnsi="Arial" w:cs="Arial"/><w:color w:val="0D0D0D" w:themeColor="text1" w:themeTint="F2"/><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr><w:t>Ser um aplicativo para celular e tablete que economize o consumo de energia do aparelho.</w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="Arial" w:hAnsi="Arial" w:cs="Arial"/><w:color w:val="0D0D0D" w:themeColor="text1" w:themeTint="F2"/><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr><w:t xml:space="preserve"> Gastando menos energia elétrica na hora de carrega-lo </w:t></w:r></w:p><w:p><w:pPr><w:rPr><w:rFonts w:ascii="Arial" w:hAnsi="Arial" w:cs="Arial"/><w:color w:val="0D0D0D" w:themeColor="text1" w:themeTint="F2"/><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr></w:pPr></w:p><w:p><w:pPr><w:rPr><w:rFonts w:ascii="Arial" w:hAnsi="Arial" w:cs="Arial"/><w:color w:val="0D0D0D" w:themeColor="text1" w:themeTint="F2"/><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr></w:pPr><w:bookmarkStart w:id="0" w:name="_GoBack"/><w:bookmarkEnd w:id="0"/></w:p><w:p><w:pPr><w:rPr><w:rFonts w:ascii="Arial" w:hAnsi="Arial" w:cs="Arial"/><w:color w:val="0D0D0D" w:themeColor="text1" w:themeTint="F2"/><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:ascii="Arial" w:hAnsi="Arial" w:cs="Arial"/><w:b/><w:color w:val="0D0D0D" w:themeColor="text1" w:themeTint="F2"/><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr><w:t>Como</w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="Arial" w:hAnsi="Arial" w:cs="Arial"/><w:color w:val="0D0D0D" w:themeColor="text1" w:themeTint="F2"/><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr><w:t>:</w:t></w:r></w:p><w:p><w:pPr><w:rPr><w:rFonts w:ascii="Arial" w:hAnsi="Arial" w:cs="Arial"/><w:color w:val="0D0D0D" w:themeColor="text1" w:themeTint="F2"/><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr></w:pPr></w:p><w:p><w:pPr><w:rPr><w:rFonts w:ascii="Arial" w:hAnsi="Arial" w:cs="Arial"/><w:color w:val="0D0D0D" w:themeColor="text1" w:themeTint="F2"/><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:ascii="Arial" w:hAnsi="Arial" w:cs="Arial"/><w:color w:val="0D0D0D" w:themeColor="text1" w:themeTint="F2"/><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr><w:t>- Fechar aplicativos de segundo plano.</w:t></w:r></w:p><w:p><w:pPr><w:rPr><w:rFonts w:ascii="Arial" w:hAnsi="Arial" w:cs="Arial"/><w:color w:val="0D0D0D" w:themeColor="text1" w:themeTint="F2"/><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:ascii="Arial" w:hAnsi="Arial" w:cs="Arial"/><w:color w:val="0D0D0D" w:themeColor="text1" w:themeTint="F2"/><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr><w:t>- Diminuir brilho da tela.</w:t></w:r></w:p><w:p><w:pPr><w:rPr><w:rFonts w:ascii="Arial" w:hAnsi="Arial" w:cs="Arial"/><w:color w:val="0D0D0D" w:themeColor="text1" w:themeTint="F2"/><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:ascii="Arial" w:hAnsi="Arial" w:cs="Arial"/><w:color w:val="0D0D0D" w:themeColor="text1" w:themeTint="F2"/><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr><w:t>- Ativar modo pouca energia.</w:t></w:r></w:p><w:p><w:pPr><w:rPr><w:rFonts w:ascii="Arial" w:hAnsi="Arial" w:cs="Arial"/><w:color w:val="0D0D0D" w:themeColor="text1" w:themeTint="F2"/><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:ascii="Arial" w:hAnsi="Arial" w:cs="Arial"/><w:color w:val="0D0D0D" w:themeColor="text1" w:themeTint="F2"/><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr><w:t>- Diminuir taxa de atualização da tela.</w:t></w:r></w:p><w:p><w:pPr><w:rPr><w:rFonts w:ascii="Arial" w:hAnsi="Arial" w:cs="Arial"/><w:color w:val="0D0D0D" w:themeColor="text1" w:themeTint="F2"/><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:ascii="Arial" w:hAnsi="Arial" w:cs="Arial"/><w:color w:val="0D0D0D" w:themeColor="text1" w:themeTint="F2"/><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr><w:t>- Relatório de uso de bateria.</w:t></w:r></w:p><w:p><w:pPr><w:rPr><w:rFonts w:ascii="Arial" w:hAnsi="Arial" w:cs="Arial"/><w:color w:val="0D0D0D" w:themeColor="text1" w:themeTint="F2"/><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:ascii="Arial" w:hAnsi="Arial" w:cs="Arial"/><w:color w:val="0D0D0D" w:themeColor="text1" w:themeTint="F2"/><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr><w:t>- Mensagens de pausa de uso.</w:t></w:r></w:p><w:p><w:pPr><w:rPr><w:rFonts w:ascii="Arial" w:hAnsi="Arial" w:cs="Arial"/><w:color w:val="0D0D0D" w:themeColor="text1" w:themeTint="F2"/><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:ascii="Arial" w:hAnsi="Arial" w:cs="Arial"/><w:color w:val="0D0D0D" w:themeColor="text1" w:themeTint="F2"/><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr><w:t>- Desativar funcionalidades que não estão em uso. (</w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:rPr><w:rFonts w:ascii="Arial" w:hAnsi="Arial" w:cs="Arial"/><w:color w:val="0D0D0D" w:themeColor="text1" w:themeTint="F2"/><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr><w:t>Wifi</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:rPr><w:rFonts w:ascii="Arial" w:hAnsi="Arial" w:cs="Arial"/><w:color w:val="0D0D0D" w:themeColor="text1" w:themeTint="F2"/><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr><w:t>, Bluetooth)</w:t></w:r></w:p><w:p><w:pPr><w:rPr><w:rFonts w:ascii="Arial" w:hAnsi="Arial" w:cs="Arial"/><w:color w:val="0D0D0D" w:themeColor="text1" w:themeTint="F2"/><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:ascii="Arial" w:hAnsi="Arial" w:cs="Arial"/><w:color w:val="0D0D0D" w:themeColor="text1" w:themeTint="F2"/><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr><w:t>- Dica do dia no aplicativo.</w:t></w:r></w:p><w:p><w:pPr><w:rPr><w:color w:val="0D0D0D" w:themeColor="text1" w:themeTint="F2"/><w:sz w:val="28"/><w:szCs w:val="28"/></w:rPr></w:pPr></w:p><w:p><w:pPr><w:rPr><w:color w:val="0D0D0D" w:themeColor="text1" w:themeTint="F2"/><w:sz w:val="28"/><w:szCs w:val="28"/></w:rPr></w:pPr></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
$d.Content.InsertXML($xml)
